$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1285113333333333
$ws.Range("H2").Value = 0.385534
$ws.Range("I2").Value = 0.03749201237720504
$ws.Range("J2").Value = 0.03749201237720504
$ws.Range("Q2").Value = 0.003459482256222222
$ws.Range("R2").Value = 0.031135340306
$ws.Range("S2").Value = 0.03749201237720504
$ws.Range("T2").Value = 0.03749201237720504

# Row 3
$ws.Range("I3").Value = 0.7552862722193517
$ws.Range("J3").Value = 0.755286272219352
$ws.Range("S3").Value = 0.7552862722193517
$ws.Range("T3").Value = 0.755286272219352

# Row 4
$ws.Range("G4").Value = 0.692415
$ws.Range("H4").Value = 2.077245
$ws.Range("I4").Value = 0.2020057770533527
$ws.Range("J4").Value = 0.2020057770533527
$ws.Range("Q4").Value = 0.018639580995
$ws.Range("R4").Value = 0.167756228955
$ws.Range("S4").Value = 0.2020057770533527
$ws.Range("T4").Value = 0.2020057770533527

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01787866666666667
$ws.Range("H5").Value = 0.053636
$ws.Range("I5").Value = 0.005215938350090445
$ws.Range("J5").Value = 0.005215938350090446
$ws.Range("Q5").Value = 0.0004812877471111111
$ws.Range("R5").Value = 0.004331589724
$ws.Range("S5").Value = 0.005215938350090445
$ws.Range("T5").Value = 0.005215938350090446
